$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.784.97"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.20%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.392.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.34%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.38"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.46"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.531"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.83%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0824"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "31.07"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -6.38%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.778.51"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.75%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.67"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.24"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.32%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.412.42"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.761"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.92%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "40.759.90"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0912"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.25%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.16"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.99"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.75"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -5.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.92"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.64"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.93%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.70"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.48%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.41"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.83"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.72"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.95%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.23"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.60%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0730"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.45"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.80"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.56%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "15.94"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -7.90%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.74"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -8.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0988"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.82"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.94%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -7.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.972.70"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0271"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.75"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.81"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.26"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.643.42"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.57%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "72.93"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "93.21"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.58"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.53%  "
